$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "super_smaller_1_first" / "super_smaller_1_second" pair of rows
# (rows 15-16), since those CSV files are no longer used.
$ws.Rows("15:16").Delete()

# Row 14 ("second" section header) no longer needs the thick-bottom border /
# explicit row height that used to line up visually with the deleted row 15.
$ws.Rows("14").AutoFit()

# "first" section (rows 4-11): resimulated, so the second run's test-file
# column now points at the "_second" file instead of repeating "_first".
$ws.Range("I5").Value = "super_1_second"
$ws.Range("I7").Value = "super_2_second"
$ws.Range("I9").Value = "super_3_second"
$ws.Range("I11").Value = "super_4_second"

# "second" section (now rows 15-23 after the row deletion above): same
# resimulation update for the remaining super_smaller_* pairs.
$ws.Range("I16").Value = "super_smaller_2_second"
$ws.Range("I18").Value = "super_smaller_3_second"
$ws.Range("I20").Value = "super_smaller_4_second"
$ws.Range("I22").Value = "super_smaller_12_second"
$ws.Range("I23").Value = "super_smaller_12_third"

# Restore the view state captured when the edit was saved.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Application.ActiveWindow.Zoom = 121
$ws.Range("D17").Select()
